$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Requisitos:" block lists three course requirement lines in B23:C25.
# Reorder them so the "LOM3246 ... (Indicação de Conjunto)" line moves from
# the first position (row 23) to the last position (row 25), while the
# other two lines shift up by one row.

$reqConjunto = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)`n"
$reqFisicaIV = "LOB1021 -  Física IV  (Requisito)`n"
$reqIntroducao = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"

$ws.Range("B23").Value = $reqFisicaIV
$ws.Range("C23").Value = $reqFisicaIV

$ws.Range("B24").Value = $reqIntroducao
$ws.Range("C24").Value = $reqIntroducao

$ws.Range("B25").Value = $reqConjunto
$ws.Range("C25").Value = $reqConjunto
